# Disable Runmode ("N") for every Test Suite row except the Profile module
# (D Suite / row 5), which is being left "Y" since new Profile-comments
# scripts are being added for it in this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

# Move the active selection from C8 to C11, as recorded in the sheet view.
[void]$ws.Range("C11").Select()
